# ---------------------------------------------------------------------------
# Add 2022-Q4 data:
#   1. Insert a new "2022-Q4" row at the top of the "总计" (summary) sheet,
#      pushing the existing quarters down by one row.
#   2. Insert a brand-new "2022-Q4" worksheet (positioned right after "总计",
#      before "2022-Q3") holding the per-fund holdings detail for the quarter.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" summary sheet — prepend the 2022-Q4 row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Snapshot the existing data rows (rows 2..6) before we start overwriting.
$oldRows = @()
for ($r = 2; $r -le 6; $r++) {
    $oldRows += ,@(
        $summary.Cells.Item($r, 2).Value2,
        $summary.Cells.Item($r, 3).Value2,
        $summary.Cells.Item($r, 4).Value2
    )
}

# Make sure row 7 exists with the same look as the other data rows (copy
# formatting down from row 6, the last currently-populated data row).
$summary.Range("A6:D6").Copy()
$summary.Range("A7:D7").PasteSpecial(-4122)

# Row 2 becomes the new 2022-Q4 entry.
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 13
$summary.Cells.Item(2, 4).Value = 1.93

# Rows 3..7 get the old rows 2..6, re-indexed sequentially (1..5) in column A.
for ($i = 0; $i -lt $oldRows.Count; $i++) {
    $destRow = 3 + $i
    $summary.Cells.Item($destRow, 1).Value = $i + 1
    $summary.Cells.Item($destRow, 2).Value = $oldRows[$i][0]
    $summary.Cells.Item($destRow, 3).Value = $oldRows[$i][1]
    $summary.Cells.Item($destRow, 4).Value = $oldRows[$i][2]
}

# ---------------------------------------------------------------------------
# 2) New "2022-Q4" worksheet with the fund-holdings detail.
# ---------------------------------------------------------------------------
# Clone the existing "2022-Q3" sheet (current position 2) so the new sheet
# inherits identical formatting/styles, then place the clone right after
# "总计" and rename it.
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($null, $summary)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template only has 2 data rows (rows 2-3); we need 13 (rows 2-14).
# Extend the formatting down by copying row 2's look into rows 4-14
# (rows 2 and 3 already carry the right per-column formatting).
$q4.Range("A2:H2").Copy()
$q4.Range("A4:H14").PasteSpecial(-4122)

$rows = @(
    @(0,  "002446", "广发利鑫灵活配置混合A",                       "22.53", "73.90", "4.05", "0.9125", 3),
    @(1,  "398021", "中海能源策略混合",                             "18.04", "90.92", "3.03", "0.5466", 10),
    @(2,  "011172", "广发利鑫灵活配置混合C",                       "7.03",  "73.90", "4.05", "0.2847", 3),
    @(3,  "000963", "兴业多策略灵活配置混合",                       "1.60",  "87.57", "3.72", "0.0595", 10),
    @(4,  "011446", "长江新能源产业混合A",                         "1.79",  "85.35", "2.79", "0.0499", 10),
    @(5,  "007251", "广发睿享稳健增利混合A",                       "1.26",  "39.64", "2.28", "0.0287", 4),
    @(6,  "970113", "兴证资管金麒麟兴睿优选一年持有期混合B",       "0.67",  "84.89", "3.44", "0.0230", 6),
    @(7,  "011447", "长江新能源产业混合C",                         "0.48",  "85.35", "2.79", "0.0134", 10),
    @(8,  "010765", "国寿安保华丰混合A",                           "0.40",  "83.84", "2.40", "0.0096", 4),
    @(9,  "010766", "国寿安保华丰混合C",                           "0.01",  "83.84", "2.40", "0.0002", 4),
    @(10, "970112", "兴证资管金麒麟兴睿优选一年持有期混合A",       "0.00",  "84.89", "3.44", $null,    6),
    @(11, "970114", "兴证资管金麒麟兴睿优选一年持有期混合C",       "0.00",  "84.89", "3.44", $null,    6),
    @(12, "011702", "广发睿享稳健增利混合C",                       "0.00",  "39.64", "2.28", $null,    4)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 2 + $i
    $row = $rows[$i]

    $q4.Cells.Item($r, 1).Value = $row[0]

    # Text-like columns: force text storage so "22.53" etc. isn't reinterpreted
    # as a number.
    $q4.Range($q4.Cells.Item($r, 2), $q4.Cells.Item($r, 7)).NumberFormat = "@"
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]

    if ($null -eq $row[6]) {
        # Last three rows store the market-value column as a real 0, not text.
        $q4.Cells.Item($r, 7).NumberFormat = "General"
        $q4.Cells.Item($r, 7).Value = 0
    } else {
        $q4.Cells.Item($r, 7).Value = $row[6]
    }

    $q4.Cells.Item($r, 8).Value = $row[7]
}

# Restore the original active-sheet state: "2021-Q1" (now the last tab)
# was the selected sheet before this edit.
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$last.Select()
$last.Range("A1").Select()
